$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new "Meta description: ..." paragraph right after the
#    Heading1 title paragraph.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaText = "Meta description: Enjoy Dragon Stone - a traditional game with 5 reels, 20 paylines and big winnings up to 698x value of bet. Exciting dragon features and abilities included."
$metaPara.Range.InsertAfter($metaText)

# Bold just the leading "Meta description" label.
$boldRng = $metaPara.Range.Duplicate
$boldRng.Find.Execute("Meta description", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$boldRng.Bold = 1

# ---------------------------------------------------------------------
# 2) Near the end of the document: drop the bold "Play Dragon Stone
#    Free..." paragraph entirely, and rewrite the italic paragraph's
#    text into the new image-generation prompt.
# ---------------------------------------------------------------------
$boldTitleText = "Play Dragon Stone Free - Big Winnings & Exciting Dragon Features"
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $boldTitleText) {
        $p.Range.Delete()
        break
    }
}

$oldPromptText = "Enjoy Dragon Stone - a traditional game with 5 reels, 20 paylines and big winnings up to 698x value of bet. Exciting dragon features and abilities included."

# Locate the italic paragraph by walking from the end of the document
# (the prompt paragraph is always the very last body paragraph), so we
# never touch the unrelated occurrence of this sentence in the new
# meta-description paragraph near the top.
$promptPara = $null
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $oldPromptText) {
        $promptPara = $p
        break
    }
}

# Placeholder tokens stand in for the straight quote / apostrophe
# characters so Word's smart-quote autocorrect (triggered by
# Find.Execute replacements) can't mangle them; we swap the tokens
# back in afterwards via plain Range.Text assignment, which bypasses
# autocorrect.
$newPromptPlaceholder = "Prompt: Create a cartoon-style feature image for the game TOK1Dragon StoneTOK2. The image should feature a happy Maya warrior with glasses. The image should have a colorful background that features the four dragons of different colors with the warrior standing in the center. The Maya warrior should be wearing a headdress and glasses with a big smile on their face, holding a dragonTOK3s stone in their hand. The dragonTOK4s stone should glow and have a bright aura surrounding it. The warriorsTOK5 clothing should be brightly colored to match the dragons around them. The image should be dynamic and engaging, representing the adventurous and fun nature of the game."

$promptRng = $promptPara.Range
$promptRng.Find.Execute($oldPromptText, $true, $false, $false, $false, $false, $true, 1, $false, $newPromptPlaceholder, 2) | Out-Null

$tokenMap = @{
    "TOK1" = '"'
    "TOK2" = '"'
    "TOK3" = "'"
    "TOK4" = "'"
    "TOK5" = "'"
}

foreach ($tok in $tokenMap.Keys) {
    $tr = $promptPara.Range.Duplicate
    $tr.Find.Execute($tok, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $tr.Text = $tokenMap[$tok]
}

Write-Output "done"
